$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.127.32"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.544.01"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'318.36"
$ws.Range("E5").Value = "  +4.16%  "
$ws.Range("D6").Value = "'95.99"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'36.45"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'7.64"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "2.937.22"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.575.60"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.55"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "'0.854"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "43.095.12"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'13.13"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'6.68"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "0.0₃0973"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'70.51"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'252.86"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "'2.03"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'27.11"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("D29").Value = "'40.14"
$ws.Range("E29").Value = "  +4.41%  "
$ws.Range("D30").Value = "'10.25"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'155.59"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "'2.14"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "'3.35"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'19.03"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("D36").Value = "'0.0793"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("D39").Value = "'24.41"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = "  +9.89%  "
$ws.Range("D42").Value = "'3.39"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "'3.83"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "'0.0305"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "2.020.07"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "'86.01"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'8.87"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "2.792.37"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'74.61"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").Value = "'102.77"
$ws.Range("E51").Value = "  -0.31%  "
